$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Se préinscrire"
$ws.Range("B2").Value = "null"
$ws.Range("C2").Value = "En cours"

$ws.Range("A3").Value = "Se préinscrire"
$ws.Range("B3").Value = "null"
$ws.Range("C3").Value = "Terminé"

$ws.Range("A4").Value = "Se préinscrire"
$ws.Range("B4").Value = "null"
$ws.Range("C4").Value = "Terminé"

$ws.Range("A5").Value = "Se préinscrire"
$ws.Range("B5").Value = "null"
$ws.Range("C5").Value = "En cours"

$ws.Range("A6").Value = "Se préinscrire"
$ws.Range("B6").Value = "null"
$ws.Range("C6").Value = "Terminé"
